$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 1301.0834  # H28: was 1222.1818
$ws.Cells.Item(28, 10).Value = 2439.6  # J28: was 2507.25
$ws.Cells.Item(28, 12).Value = 2439.6  # L28: was 2507.25
$ws.Cells.Item(28, 14).Value = -3409.6  # N28: was -3477.25
# Row 62
$ws.Cells.Item(62, 8).Value = 3920  # H62: was 4037.5
$ws.Cells.Item(62, 9).Value = 3566.6667  # I62: was 3625
$ws.Cells.Item(62, 11).Value = 3566.6667  # K62: was 3625
$ws.Cells.Item(62, 13).Value = -2942.6667  # M62: was -3001
# Row 65
$ws.Cells.Item(65, 8).Value = 3920  # H65: was 4037.5
$ws.Cells.Item(65, 9).Value = 3566.6667  # I65: was 3625
$ws.Cells.Item(65, 11).Value = 17833.3335  # K65: was 18125
$ws.Cells.Item(65, 13).Value = -14713.3335  # M65: was -15005
# Row 92
$ws.Cells.Item(92, 8).Value = 47619372  # H92: was 41666996
$ws.Cells.Item(92, 9).Value = 55555864  # I92: was 55555876
$ws.Cells.Item(92, 10).Value = 434  # J92: was 358.5
$ws.Cells.Item(92, 11).Value = 55555864  # K92: was 55555876
$ws.Cells.Item(92, 12).Value = 434  # L92: was 358.5
$ws.Cells.Item(92, 13).Value = -55554616  # M92: was -55554628
$ws.Cells.Item(92, 14).Value = -2930  # N92: was -2854.5
# Row 107
$ws.Cells.Item(107, 8).Value = 75779.336  # H107: was 69991.46000000001
$ws.Cells.Item(107, 9).Value = 75779.336  # I107: was 69991.46000000001
$ws.Cells.Item(107, 11).Value = 75779.336  # K107: was 69991.46000000001
$ws.Cells.Item(107, 13).Value = -73859.336  # M107: was -68071.46000000001
# Row 138
$ws.Cells.Item(138, 8).Value = 3196.8667  # H138: was 2697.7
$ws.Cells.Item(138, 9).Value = 1132  # I138: was 897.6667
$ws.Cells.Item(138, 10).Value = 3426.2964  # J138: was 3469.1428
$ws.Cells.Item(138, 11).Value = 3396  # K138: was 2693.0001
$ws.Cells.Item(138, 12).Value = 10278.8892  # L138: was 10407.4284
$ws.Cells.Item(138, 13).Value = 1744  # M138: was 2446.9999
$ws.Cells.Item(138, 14).Value = -20558.8892  # N138: was -20687.4284

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 41
$ws.Cells.Item(41, 8).Value = 565.2  # H41: was 587.2
$ws.Cells.Item(41, 9).Value = 565.2  # I41: was 587.2
$ws.Cells.Item(41, 11).Value = 565.2  # K41: was 587.2
$ws.Cells.Item(41, 13).Value = -151.2  # M41: was -173.2
# Row 61
$ws.Cells.Item(61, 8).Value = 2900  # H61: was 1459
$ws.Cells.Item(61, 9).Value = 2000  # I61: was 873.75
$ws.Cells.Item(61, 11).Value = 2000  # K61: was 873.75
$ws.Cells.Item(61, 13).Value = -1788  # M61: was -661.75
# Row 74
$ws.Cells.Item(74, 8).Value = 2700  # H74: was 2549
$ws.Cells.Item(74, 9).Value = 1600  # I74: was 1961.25
$ws.Cells.Item(74, 11).Value = 1600  # K74: was 1961.25
$ws.Cells.Item(74, 13).Value = -726  # M74: was -1087.25
# Row 77
$ws.Cells.Item(77, 8).Value = 2700  # H77: was 2549
$ws.Cells.Item(77, 9).Value = 1600  # I77: was 1961.25
$ws.Cells.Item(77, 11).Value = 8000  # K77: was 9806.25
$ws.Cells.Item(77, 13).Value = -3632  # M77: was -5438.25
# Row 119
$ws.Cells.Item(119, 8).Value = 50000  # H119: was 47500
$ws.Cells.Item(119, 10).Value = 50000  # J119: was 47500
$ws.Cells.Item(119, 12).Value = 50000  # L119: was 47500
$ws.Cells.Item(119, 14).Value = -59676  # N119: was -57176
# Row 122
$ws.Cells.Item(122, 8).Value = 10000  # H122: was 5179.6
$ws.Cells.Item(122, 9).Value = 10000  # I122: was 5999.6665
$ws.Cells.Item(122, 10).Value = 0  # J122: was 3949.5
$ws.Cells.Item(122, 11).Value = 30000  # K122: was 17998.9995
$ws.Cells.Item(122, 12).Value = 0  # L122: was 11848.5
$ws.Cells.Item(122, 13).Value = -27550  # M122: was -15548.9995
$ws.Cells.Item(122, 14).ClearContents()  # N122: was -16748.5
# Row 132
$ws.Cells.Item(132, 8).Value = 4059  # H132: was 3343.1667
$ws.Cells.Item(132, 9).Value = 4059  # I132: was 3343.1667
$ws.Cells.Item(132, 11).Value = 12177  # K132: was 10029.5001
$ws.Cells.Item(132, 13).Value = -9647  # M132: was -7499.500100000001
# Row 136
$ws.Cells.Item(136, 8).Value = 2900  # H136: was 1459
$ws.Cells.Item(136, 9).Value = 2000  # I136: was 873.75
$ws.Cells.Item(136, 11).Value = 6000  # K136: was 2621.25
$ws.Cells.Item(136, 13).Value = -3450  # M136: was -71.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 10
$ws.Cells.Item(10, 8).Value = 312.85715  # H10: was 415.77777
$ws.Cells.Item(10, 9).Value = 750  # I10: was 478.85715
$ws.Cells.Item(10, 10).Value = 138  # J10: was 195
$ws.Cells.Item(10, 11).Value = 750  # K10: was 478.85715
$ws.Cells.Item(10, 12).Value = 138  # L10: was 195
$ws.Cells.Item(10, 13).Value = -610  # M10: was -338.85715
$ws.Cells.Item(10, 14).Value = -418  # N10: was -475
# Row 107
$ws.Cells.Item(107, 8).Value = 45033.332  # H107: was 40523.7
$ws.Cells.Item(107, 9).Value = 57546.43  # I107: was 57507
$ws.Cells.Item(107, 10).Value = 1237.5  # J107: was 896
$ws.Cells.Item(107, 11).Value = 57546.43  # K107: was 57507
$ws.Cells.Item(107, 12).Value = 1237.5  # L107: was 896
$ws.Cells.Item(107, 13).Value = -55626.43  # M107: was -55587
$ws.Cells.Item(107, 14).Value = -5077.5  # N107: was -4736
# Row 134
$ws.Cells.Item(134, 8).Value = 903.6667  # H134: was 904
$ws.Cells.Item(134, 9).Value = 903.6667  # I134: was 904
$ws.Cells.Item(134, 11).Value = 2711.0001  # K134: was 2712
$ws.Cells.Item(134, 13).Value = -176.0001000000002  # M134: was -177

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Cells.Item(6, 8).Value = 5716593  # H6: was 5550745.5
$ws.Cells.Item(6, 9).Value = 8001230  # I6: was 5550745.5
$ws.Cells.Item(6, 10).Value = 5000  # J6: was 0
$ws.Cells.Item(6, 11).Value = 8001230  # K6: was 5550745.5
$ws.Cells.Item(6, 12).Value = 5000  # L6: was 0
$ws.Cells.Item(6, 13).Value = -8001117  # M6: was -5550632.5
$ws.Cells.Item(6, 14).Value = -5226  # N6: was None
# Row 10
$ws.Cells.Item(10, 8).Value = 3984.5  # H10: was 855
$ws.Cells.Item(10, 9).Value = 470  # I10: was 32.5
$ws.Cells.Item(10, 10).Value = 7499  # J10: was 2500
$ws.Cells.Item(10, 11).Value = 470  # K10: was 32.5
$ws.Cells.Item(10, 12).Value = 7499  # L10: was 2500
$ws.Cells.Item(10, 13).Value = -331  # M10: was 106.5
$ws.Cells.Item(10, 14).Value = -7777  # N10: was -2778
# Row 16
$ws.Cells.Item(16, 8).Value = 532  # H16: was 448.75
$ws.Cells.Item(16, 9).Value = 698  # I16: was 365.66666
$ws.Cells.Item(16, 10).Value = 449  # J16: was 698
$ws.Cells.Item(16, 11).Value = 698  # K16: was 365.66666
$ws.Cells.Item(16, 12).Value = 449  # L16: was 698
$ws.Cells.Item(16, 13).Value = -411  # M16: was -78.66665999999998
$ws.Cells.Item(16, 14).Value = -1023  # N16: was -1272
# Row 22
$ws.Cells.Item(22, 8).Value = 999  # H22: was 997.75
$ws.Cells.Item(22, 9).Value = 999  # I22: was 997.75
$ws.Cells.Item(22, 11).Value = 999  # K22: was 997.75
$ws.Cells.Item(22, 13).Value = -649  # M22: was -647.75
# Row 58
$ws.Cells.Item(58, 8).Value = 1882.6666  # H58: was 2011.5
$ws.Cells.Item(58, 9).Value = 1859.2  # I58: was 2011.5
$ws.Cells.Item(58, 10).Value = 2000  # J58: was 0
$ws.Cells.Item(58, 11).Value = 1859.2  # K58: was 2011.5
$ws.Cells.Item(58, 12).Value = 2000  # L58: was 0
$ws.Cells.Item(58, 13).Value = -1656.2  # M58: was -1808.5
$ws.Cells.Item(58, 14).Value = -2406  # N58: was None
# Row 103
$ws.Cells.Item(103, 8).Value = 37485.5  # H103: was 39981
$ws.Cells.Item(103, 9).Value = 37485.5  # I103: was 39981
$ws.Cells.Item(103, 11).Value = 37485.5  # K103: was 39981
$ws.Cells.Item(103, 13).Value = -36313.5  # M103: was -38809
# Row 113
$ws.Cells.Item(113, 8).Value = 532  # H113: was 448.75
$ws.Cells.Item(113, 9).Value = 698  # I113: was 365.66666
$ws.Cells.Item(113, 10).Value = 449  # J113: was 698
$ws.Cells.Item(113, 11).Value = 698  # K113: was 365.66666
$ws.Cells.Item(113, 12).Value = 449  # L113: was 698
$ws.Cells.Item(113, 13).Value = 1472  # M113: was 1804.33334
$ws.Cells.Item(113, 14).Value = -4789  # N113: was -5038
# Row 132
$ws.Cells.Item(132, 8).Value = 1222.625  # H132: was 1195.6666
$ws.Cells.Item(132, 9).Value = 963.5  # I132: was 965.8570999999999
$ws.Cells.Item(132, 11).Value = 2890.5  # K132: was 2897.5713
$ws.Cells.Item(132, 13).Value = -360.5  # M132: was -367.5712999999996
# Row 134
$ws.Cells.Item(134, 8).Value = 1441.0526  # H134: was 1378.9048
$ws.Cells.Item(134, 9).Value = 1211.25  # I134: was 1164.2778
$ws.Cells.Item(134, 11).Value = 3633.75  # K134: was 3492.8334
$ws.Cells.Item(134, 13).Value = -1098.75  # M134: was -957.8334000000004
# Row 136
$ws.Cells.Item(136, 8).Value = 1882.6666  # H136: was 2011.5
$ws.Cells.Item(136, 9).Value = 1859.2  # I136: was 2011.5
$ws.Cells.Item(136, 10).Value = 2000  # J136: was 0
$ws.Cells.Item(136, 11).Value = 5577.6  # K136: was 6034.5
$ws.Cells.Item(136, 12).Value = 6000  # L136: was 0
$ws.Cells.Item(136, 13).Value = -3027.6  # M136: was -3484.5
$ws.Cells.Item(136, 14).Value = -11100  # N136: was None
# Row 141
$ws.Cells.Item(141, 8).Value = 761109.7  # H141: was 883331.6
$ws.Cells.Item(141, 10).Value = 761109.7  # J141: was 883331.6
$ws.Cells.Item(141, 12).Value = 761109.7  # L141: was 883331.6
$ws.Cells.Item(141, 14).Value = -771469.7  # N141: was -893691.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Cells.Item(63, 8).Value = 2341.3333  # H63: was 3006
# Row 66
$ws.Cells.Item(66, 8).Value = 2341.3333  # H66: was 3006
# Row 114
$ws.Cells.Item(114, 8).Value = 1121.1666  # H114: was 1101.9231
$ws.Cells.Item(114, 10).Value = 1254.4286  # J114: was 1206.5
$ws.Cells.Item(114, 12).Value = 3763.2858  # L114: was 3619.5
$ws.Cells.Item(114, 14).Value = -10271.2858  # N114: was -10127.5
# Row 117
$ws.Cells.Item(117, 8).Value = 3340.9092  # H117: was 3670
$ws.Cells.Item(117, 9).Value = 114.5  # I117: was 117
$ws.Cells.Item(117, 10).Value = 4057.889  # J117: was 4558.25
$ws.Cells.Item(117, 11).Value = 343.5  # K117: was 351
$ws.Cells.Item(117, 12).Value = 12173.667  # L117: was 13674.75
$ws.Cells.Item(117, 13).Value = 3098.5  # M117: was 3091
$ws.Cells.Item(117, 14).Value = -19057.667  # N117: was -20558.75
# Row 126
$ws.Cells.Item(126, 8).Value = 8000  # H126: was 12000
$ws.Cells.Item(126, 9).Value = 8000  # I126: was 12000
$ws.Cells.Item(126, 11).Value = 24000  # K126: was 36000
$ws.Cells.Item(126, 13).Value = -19060  # M126: was -31060

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Cells.Item(113, 8).Value = 923.5  # H113: was 916.8
$ws.Cells.Item(113, 9).Value = 923.5  # I113: was 916.8
$ws.Cells.Item(113, 11).Value = 923.5  # K113: was 916.8
$ws.Cells.Item(113, 13).Value = 1246.5  # M113: was 1253.2
# Row 121
$ws.Cells.Item(121, 8).Value = 20000  # H121: was 0
$ws.Cells.Item(121, 10).Value = 20000  # J121: was 0
$ws.Cells.Item(121, 12).Value = 20000  # L121: was 0
$ws.Cells.Item(121, 14).Value = -23494  # N121: was None
# Row 132
$ws.Cells.Item(132, 8).Value = 1348.75  # H132: was 1279.4
$ws.Cells.Item(132, 9).Value = 1348.75  # I132: was 1279.4
$ws.Cells.Item(132, 11).Value = 4046.25  # K132: was 3838.2
$ws.Cells.Item(132, 13).Value = -1516.25  # M132: was -1308.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Cells.Item(55, 8).Value = 1027.0435  # H55: was 1067.1818
$ws.Cells.Item(55, 9).Value = 781.5714  # I55: was 830.61536
$ws.Cells.Item(55, 11).Value = 781.5714  # K55: was 830.61536
$ws.Cells.Item(55, 13).Value = -608.5714  # M55: was -657.61536
# Row 82
$ws.Cells.Item(82, 8).Value = 4268.2856  # H82: was 4563.1665
$ws.Cells.Item(82, 9).Value = 3495.6  # I82: was 3744.75
$ws.Cells.Item(82, 11).Value = 3495.6  # K82: was 3744.75
$ws.Cells.Item(82, 13).Value = -3134.6  # M82: was -3383.75
# Row 85
$ws.Cells.Item(85, 8).Value = 4268.2856  # H85: was 4563.1665
$ws.Cells.Item(85, 9).Value = 3495.6  # I85: was 3744.75
$ws.Cells.Item(85, 11).Value = 3495.6  # K85: was 3744.75
$ws.Cells.Item(85, 13).Value = -2247.6  # M85: was -2496.75
# Row 132
$ws.Cells.Item(132, 8).Value = 1496.6666  # H132: was 1943.3334
$ws.Cells.Item(132, 9).Value = 1496.6666  # I132: was 1943.3334
$ws.Cells.Item(132, 11).Value = 4489.9998  # K132: was 5830.0002
$ws.Cells.Item(132, 13).Value = -1959.9998  # M132: was -3300.0002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 2824.4546  # H81: was 2365.4443
$ws.Cells.Item(81, 10).Value = 4917.25  # J81: was 4944.5
$ws.Cells.Item(81, 12).Value = 9834.5  # L81: was 9889
$ws.Cells.Item(81, 14).Value = -11956.5  # N81: was -12011
# Row 84
$ws.Cells.Item(84, 8).Value = 2824.4546  # H84: was 2365.4443
$ws.Cells.Item(84, 10).Value = 4917.25  # J84: was 4944.5
$ws.Cells.Item(84, 12).Value = 49172.5  # L84: was 49445
$ws.Cells.Item(84, 14).Value = -59780.5  # N84: was -60053
# Row 96
$ws.Cells.Item(96, 8).Value = 3800  # H96: was 3333.8333
$ws.Cells.Item(96, 9).Value = 4500  # I96: was 3625.75
$ws.Cells.Item(96, 11).Value = 4500  # K96: was 3625.75
$ws.Cells.Item(96, 13).Value = -3127  # M96: was -2252.75
# Row 113
$ws.Cells.Item(113, 8).Value = 469.125  # H113: was 492.73334
$ws.Cells.Item(113, 9).Value = 254.88889  # I113: was 272.375
$ws.Cells.Item(113, 11).Value = 764.6666700000001  # K113: was 817.125
$ws.Cells.Item(113, 13).Value = 1405.33333  # M113: was 1352.875
# Row 132
$ws.Cells.Item(132, 8).Value = 901  # H132: was 899.6667
$ws.Cells.Item(132, 9).Value = 901  # I132: was 899.6667
$ws.Cells.Item(132, 11).Value = 2703  # K132: was 2699.0001
$ws.Cells.Item(132, 13).Value = -173  # M132: was -169.0001000000002
# Row 136
$ws.Cells.Item(136, 8).Value = 1299.8182  # H136: was 1414.8
$ws.Cells.Item(136, 9).Value = 1299.8182  # I136: was 1414.8
$ws.Cells.Item(136, 11).Value = 3899.4546  # K136: was 4244.4
$ws.Cells.Item(136, 13).Value = -1349.4546  # M136: was -1694.4
